{"js": "// 1) Fix the typo \"hopping\" -> \"shopping\" in the intro paragraph\n//    (\"...uncover patterns in hopping habits...\" -> \"...patterns in shopping habits...\").\n//    Search narrowly on \" hopping\" (leading space) so we don't also match the\n//    \"hopping\" substring inside the existing word \"shopping\" later in the paragraph.\nconst typoResults = context.document.body.search(\" hopping\", { matchCase: true });\ntypoResults.load(\"items\");\nawait context.sync();\n\nif (typoResults.items.length > 0) {\n  typoResults.items[0].insertText(\" shopping\", \"Replace\");\n  await context.sync();\n}\n\n// 2) Append the new Q&A paragraphs (with blank-line spacers) after the\n//    existing \"Rough breakdown of tasks:\" paragraph, at the end of the body.\nconst body = context.document.body;\n\nbody.insertParagraph(\"\", \"End\");\nbody.insertParagraph(\"Sergyo: Do people shop more online when its raining?\", \"End\");\n\nbody.insertParagraph(\"\", \"End\");\nbody.insertParagraph(\n  \"Larry: Do people shop more online when the temperatures are extreme (<20)(>85)\",\n  \"End\"\n);\n\nbody.insertParagraph(\"\", \"End\");\nbody.insertParagraph(\"Emily: What season has the highest shopping?\", \"End\");\n\nbody.insertParagraph(\"\", \"End\");\nbody.insertParagraph(\"Matthew; Did the season with the most shopping have the most rain?\", \"End\");\n\nbody.insertParagraph(\"\", \"End\");\nbody.insertParagraph(\n  \"Rebecca: Did the season with the most shopping have the most extreme temperatures?\",\n  \"End\"\n);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Fix the typo \"hopping\" -> \"shopping\" in the intro paragraph\n#    (\"...uncover patterns in hopping habits...\" -> \"...patterns in shopping habits...\").\n#    Find.Execute locates the FIRST occurrence of \"hopping\" in the document,\n#    which is this one (the later \"shopping\" elsewhere doesn't match the\n#    search term \"hopping\" as Find returns the match span itself).\n$find = $d.Content\n$found = $find.Find.Execute(\"hopping\", $true, $false, $false, $false, $false, $true, 1, $false, \"\", 0)\nif ($found) {\n    $find.Collapse(1)\n    $find.InsertBefore(\"s\")\n}\n\n# 2) Append the new Q&A paragraphs (with blank-line spacers) after the\n#    existing \"Rough breakdown of tasks:\" paragraph, at the end of the body.\nfunction Add-BlankLineThenParagraph([string]$text) {\n    $tail = $d.Content\n    $tail.Collapse(0)\n    $tail.InsertParagraphAfter()\n\n    $tail2 = $d.Content\n    $tail2.Collapse(0)\n    $tail2.InsertParagraphAfter()\n\n    $tail3 = $d.Content\n    $tail3.Collapse(0)\n    $tail3.InsertAfter($text)\n}\n\nAdd-BlankLineThenParagraph(\"Sergyo: Do people shop more online when its raining?\")\nAdd-BlankLineThenParagraph(\"Larry: Do people shop more online when the temperatures are extreme (<20)(>85)\")\nAdd-BlankLineThenParagraph(\"Emily: What season has the highest shopping?\")\nAdd-BlankLineThenParagraph(\"Matthew; Did the season with the most shopping have the most rain?\")\nAdd-BlankLineThenParagraph(\"Rebecca: Did the season with the most shopping have the most extreme temperatures?\")\n"}
